$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.101.25"
$ws.Range("E2").Value = "  +5.15%  "

$ws.Range("D3").Value = "3.242.40"
$ws.Range("E3").Value = "  +2.11%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'394.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "'108.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("D7").Value = "'0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.04%  "

$ws.Range("D8").Value = "3.236.58"
$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").Value = "'38.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").Value = "'0.0996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +13.00%  "

$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "3.760.19"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").Value = "'8.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").Value = "'19.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "3.258.49"
$ws.Range("E17").Value = "  +2.61%  "

$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").Value = "'10.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").Value = "56.962.02"
$ws.Range("E20").Value = "  +4.94%  "

$ws.Range("E21").Value = "  +0.88%  "

$ws.Range("D22").Value = "'0.0000111"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.29%  "

$ws.Range("D23").Value = "'12.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "'294.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.91%  "

$ws.Range("D25").Value = "'74.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.58%  "

$ws.Range("D26").Value = "'3.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.54%  "

$ws.Range("D27").Value = "'27.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").Value = "'7.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.94%  "

$ws.Range("D29").Value = "'7.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "'11.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.03%  "

$ws.Range("E33").Value = "  -3.73%  "

$ws.Range("D34").Value = "'39.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.71%  "

$ws.Range("D35").Value = "'0.0480"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("E36").Value = "  +2.35%  "

$ws.Range("D37").Value = "'51.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("E39").Value = "  -5.38%  "

$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("D41").Value = "'134.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.66%  "

$ws.Range("E42").Value = "  +4.13%  "

$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.23%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'16.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.22%  "

$ws.Range("D46").Value = "'0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("D47").Value = "'22.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D49").Value = "2.148.35"
$ws.Range("E49").Value = "  +2.94%  "

$ws.Range("D50").Value = "'2.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.89%  "

$ws.Range("E51").Value = "  +15.08%  "
